$wb = $excel.ActiveWorkbook

# --- Update the "Yearly" sheet (2017 / June row, row 8) ---
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsYearly.Range("L8").Value = 63.16
$wsYearly.Range("M8").Value = 37.27
$wsYearly.Range("N8").Value = 24.55

# Move the saved selection/active cell on the Yearly sheet
$wsYearly.Range("D28").Select()

# --- Update the "All Time" sheet view (selection + scroll position) ---
$wsAllTime = $wb.Worksheets.Item("All Time")
$wsAllTime.Activate()
$excel.ActiveWindow.ScrollRow = 31
$wsAllTime.Range("N43").Select()
